# "Pais" worksheet: COVID-19 case counts per country, sorted descending by
# total cases (column B). This refresh updates the case counters for the
# countries whose figures moved, and - because a handful of countries
# leap-frogged their immediate neighbour in the sort order - swaps the
# country names between the affected row pairs so the table stays sorted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Octubre de 2020 a las 11:35"

# Update country names (column A) where the sorted ranking order changed
$ws.Cells.Item(21, 1).Value = "Indonesia"
$ws.Cells.Item(22, 1).Value = "Filipinas"
$ws.Cells.Item(34, 1).Value = "Polonia"
$ws.Cells.Item(35, 1).Value = "Ecuador"
$ws.Cells.Item(62, 1).Value = "Austria"
$ws.Cells.Item(63, 1).Value = "Armenia"
$ws.Cells.Item(87, 1).Value = "Eslovaquia"
$ws.Cells.Item(88, 1).Value = "Grecia"
$ws.Cells.Item(89, 1).Value = "Croacia"
$ws.Cells.Item(219, 1).Value = "Islas Salomon"
$ws.Cells.Item(220, 1).Value = "Anguila"

# Update numeric data columns B:H for rows with refreshed statistics
# Row 5
$ws.Cells.Item(5, 2).Value = 7309164
$ws.Cells.Item(5, 3).Value = 4094
$ws.Cells.Item(5, 4).Value = 6383441
$ws.Cells.Item(5, 5).Value = 814386
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 26
$ws.Cells.Item(5, 8).Value = 111337
# Row 19
$ws.Cells.Item(19, 2).Value = 384559
$ws.Cells.Item(19, 3).Value = 1600
$ws.Cells.Item(19, 4).Value = 299229
$ws.Cells.Item(19, 5).Value = 79722
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 15
$ws.Cells.Item(19, 8).Value = 5608
# Row 21
$ws.Cells.Item(21, 2).Value = 349160
$ws.Cells.Item(21, 3).Value = 4411
$ws.Cells.Item(21, 4).Value = 273661
$ws.Cells.Item(21, 5).Value = 63231
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 112
$ws.Cells.Item(21, 8).Value = 12268
# Row 22
$ws.Cells.Item(22, 2).Value = 348698
$ws.Cells.Item(22, 3).Value = 2261
$ws.Cells.Item(22, 4).Value = 294161
$ws.Cells.Item(22, 5).Value = 48040
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 50
$ws.Cells.Item(22, 8).Value = 6497
# Row 27
$ws.Cells.Item(27, 2).Value = 299502
$ws.Cells.Item(27, 3).Value = 1002
$ws.Cells.Item(27, 4).Value = 253597
$ws.Cells.Item(27, 5).Value = 43796
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 11
$ws.Cells.Item(27, 8).Value = 2109
# Row 34
$ws.Cells.Item(34, 2).Value = 149903
$ws.Cells.Item(34, 3).Value = 8099
$ws.Cells.Item(34, 4).Value = 85588
$ws.Cells.Item(34, 5).Value = 61007
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 91
$ws.Cells.Item(34, 8).Value = 3308
# Row 35
$ws.Cells.Item(35, 2).Value = 149083
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 128134
$ws.Cells.Item(35, 5).Value = 8685
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 12264
# Row 62
$ws.Cells.Item(62, 2).Value = 60224
$ws.Cells.Item(62, 3).Value = 1552
$ws.Cells.Item(62, 4).Value = 46798
$ws.Cells.Item(62, 5).Value = 12549
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 5
$ws.Cells.Item(62, 8).Value = 877
# Row 63
$ws.Cells.Item(63, 2).Value = 59995
$ws.Cells.Item(63, 3).Value = 1371
$ws.Cells.Item(63, 4).Value = 47119
$ws.Cells.Item(63, 5).Value = 11830
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 7
$ws.Cells.Item(63, 8).Value = 1046
# Row 87
$ws.Cells.Item(87, 2).Value = 24225
$ws.Cells.Item(87, 3).Value = 1929
$ws.Cells.Item(87, 4).Value = 6926
$ws.Cells.Item(87, 5).Value = 17228
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 5
$ws.Cells.Item(87, 8).Value = 71
# Row 88
$ws.Cells.Item(88, 2).Value = 23495
$ws.Cells.Item(88, 3).Value = 0
$ws.Cells.Item(88, 4).Value = 9989
$ws.Cells.Item(88, 5).Value = 13037
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 469
# Row 89
$ws.Cells.Item(89, 2).Value = 22534
$ws.Cells.Item(89, 3).Value = 793
$ws.Cells.Item(89, 4).Value = 18628
$ws.Cells.Item(89, 5).Value = 3562
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 10
$ws.Cells.Item(89, 8).Value = 344
# Row 102
$ws.Cells.Item(102, 2).Value = 12944
$ws.Cells.Item(102, 3).Value = 241
$ws.Cells.Item(102, 4).Value = 9100
$ws.Cells.Item(102, 5).Value = 3494
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 350
# Row 104
$ws.Cells.Item(104, 2).Value = 11255
$ws.Cells.Item(104, 3).Value = 67
$ws.Cells.Item(104, 4).Value = 10360
$ws.Cells.Item(104, 5).Value = 825
$ws.Cells.Item(104, 6).Value = 0
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 70
# Row 151
$ws.Cells.Item(151, 2).Value = 3056
$ws.Cells.Item(151, 3).Value = 114
$ws.Cells.Item(151, 4).Value = 1329
$ws.Cells.Item(151, 5).Value = 1686
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 41
# Row 219
$ws.Cells.Item(219, 2).Value = 3
$ws.Cells.Item(219, 3).Value = 1
$ws.Cells.Item(219, 4).Value = 0
$ws.Cells.Item(219, 5).Value = 3
$ws.Cells.Item(219, 6).Value = 0
$ws.Cells.Item(219, 7).Value = 0
$ws.Cells.Item(219, 8).Value = 0
# Row 220
$ws.Cells.Item(220, 2).Value = 3
$ws.Cells.Item(220, 3).Value = 0
$ws.Cells.Item(220, 4).Value = 3
$ws.Cells.Item(220, 5).Value = 0
$ws.Cells.Item(220, 6).Value = 0
$ws.Cells.Item(220, 7).Value = 0
$ws.Cells.Item(220, 8).Value = 0
